$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 540.36365
$ws.Range("I33").Value = 424.25
$ws.Range("J33").Value = 850
$ws.Range("K33").Value = 424.25
$ws.Range("L33").Value = 850
$ws.Range("M33").Value = -195.25
$ws.Range("N33").Value = -1308

$ws.Range("H100").Value = 10000
$ws.Range("I100").Value = 10000
$ws.Range("K100").Value = 10000
$ws.Range("M100").Value = -9459

$ws.Range("H112").Value = 2402.8333
$ws.Range("J112").Value = 2752.0667
$ws.Range("L112").Value = 8256.2001
$ws.Range("N112").Value = -10472.2001

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 21172.139
$ws.Range("I132").Value = 21764.2
$ws.Range("J132").Value = 450
$ws.Range("K132").Value = 65292.60000000001
$ws.Range("L132").Value = 1350
$ws.Range("M132").Value = -62762.60000000001
$ws.Range("N132").Value = -6410

$ws.Range("H141").Value = 1042.75
$ws.Range("I141").Value = 906
$ws.Range("K141").Value = 2718
$ws.Range("M141").Value = 2462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18126.69
$ws.Range("I32").Value = 19277.824
$ws.Range("K32").Value = 19277.824
$ws.Range("M32").Value = -18990.824

$ws.Range("H61").Value = 5627.4585
$ws.Range("I61").Value = 970.7895
$ws.Range("K61").Value = 970.7895
$ws.Range("M61").Value = -758.7895

$ws.Range("H74").Value = 245027
$ws.Range("I74").Value = 375650.06
$ws.Range("J74").Value = 12808.223
$ws.Range("K74").Value = 375650.06
$ws.Range("L74").Value = 12808.223
$ws.Range("M74").Value = -374776.06
$ws.Range("N74").Value = -14556.223

$ws.Range("H77").Value = 245027
$ws.Range("I77").Value = 375650.06
$ws.Range("J77").Value = 12808.223
$ws.Range("K77").Value = 1878250.3
$ws.Range("L77").Value = 64041.115
$ws.Range("M77").Value = -1873882.3
$ws.Range("N77").Value = -72777.11499999999

$ws.Range("H88").Value = 6876.9
$ws.Range("I88").Value = 2092.5
$ws.Range("J88").Value = 10066.5
$ws.Range("K88").Value = 2092.5
$ws.Range("L88").Value = 10066.5
$ws.Range("M88").Value = -1686.5
$ws.Range("N88").Value = -10878.5

$ws.Range("H91").Value = 6876.9
$ws.Range("I91").Value = 2092.5
$ws.Range("J91").Value = 10066.5
$ws.Range("K91").Value = 2092.5
$ws.Range("L91").Value = 10066.5
$ws.Range("M91").Value = -688.5
$ws.Range("N91").Value = -12874.5

$ws.Range("H97").Value = 1345.7391
$ws.Range("I97").Value = 1103.05
$ws.Range("K97").Value = 1103.05
$ws.Range("M97").Value = -607.05

$ws.Range("H122").Value = 2891.6904
$ws.Range("I122").Value = 2795.1936
$ws.Range("J122").Value = 3163.6365
$ws.Range("K122").Value = 8385.5808
$ws.Range("L122").Value = 9490.9095
$ws.Range("M122").Value = -5935.5808
$ws.Range("N122").Value = -14390.9095

$ws.Range("H132").Value = 1783.3846
$ws.Range("I132").Value = 1229.4286
$ws.Range("K132").Value = 3688.2858
$ws.Range("M132").Value = -1158.2858

$ws.Range("H136").Value = 5627.4585
$ws.Range("I136").Value = 970.7895
$ws.Range("K136").Value = 2912.3685
$ws.Range("M136").Value = -362.3685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 16890.953
$ws.Range("I20").Value = 27854.584
$ws.Range("J20").Value = 2272.7778
$ws.Range("K20").Value = 27854.584
$ws.Range("L20").Value = 2272.7778
$ws.Range("M20").Value = -27607.584
$ws.Range("N20").Value = -2766.7778

$ws.Range("H99").Value = 3050.4666
$ws.Range("I99").Value = 1641.2222
$ws.Range("K99").Value = 1641.2222
$ws.Range("M99").Value = -143.2221999999999

$ws.Range("H107").Value = 3892.7837
$ws.Range("I107").Value = 2994.7036
$ws.Range("K107").Value = 2994.7036
$ws.Range("M107").Value = -1074.7036

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 113839.89
$ws.Range("I132").Value = 201061.8
$ws.Range("J132").Value = 4812.5
$ws.Range("K132").Value = 603185.3999999999
$ws.Range("L132").Value = 14437.5
$ws.Range("M132").Value = -600655.3999999999
$ws.Range("N132").Value = -19497.5

$ws.Range("H134").Value = 2292.4375
$ws.Range("I134").Value = 2128.5186
$ws.Range("K134").Value = 6385.5558
$ws.Range("M134").Value = -3850.5558

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 125000530
$ws.Range("I9").Value = 600
$ws.Range("K9").Value = 1800
$ws.Range("M9").Value = -1576

$ws.Range("H22").Value = 4759.222
$ws.Range("I22").Value = 2749.5
$ws.Range("J22").Value = 4920
$ws.Range("K22").Value = 8248.5
$ws.Range("L22").Value = 14760
$ws.Range("M22").Value = -8079.5
$ws.Range("N22").Value = -15098

$ws.Range("H23").Value = 401.06668
$ws.Range("I23").Value = 193.6
$ws.Range("K23").Value = 580.8
$ws.Range("M23").Value = -345.8

$ws.Range("H27").Value = 4759.222
$ws.Range("I27").Value = 2749.5
$ws.Range("J27").Value = 4920
$ws.Range("K27").Value = 8248.5
$ws.Range("L27").Value = 14760
$ws.Range("M27").Value = -8146.5
$ws.Range("N27").Value = -14964

$ws.Range("H139").Value = 32515
$ws.Range("I139").Value = 32515
$ws.Range("K139").Value = 97545
$ws.Range("M139").Value = -92405

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6357.2
$ws.Range("I70").Value = 5399.2
$ws.Range("K70").Value = 5399.2
$ws.Range("M70").Value = -5129.2

$ws.Range("H73").Value = 6357.2
$ws.Range("I73").Value = 5399.2
$ws.Range("K73").Value = 5399.2
$ws.Range("M73").Value = -4463.2

$ws.Range("H119").Value = 100000
$ws.Range("J119").Value = 100000
$ws.Range("L119").Value = 100000
$ws.Range("N119").Value = -109676

$ws.Range("H123").Value = 43076.92
$ws.Range("J123").Value = 43076.92
$ws.Range("L123").Value = 43076.92
$ws.Range("N123").Value = -47976.92

$ws.Range("H132").Value = 2214.4866
$ws.Range("I132").Value = 2070.2188
$ws.Range("K132").Value = 6210.6564
$ws.Range("M132").Value = -3680.6564

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2296.75
$ws.Range("I40").Value = 2296.75
$ws.Range("K40").Value = 2296.75
$ws.Range("M40").Value = -2160.75

$ws.Range("H55").Value = 1821.75
$ws.Range("I55").Value = 481
$ws.Range("J55").Value = 2626.2
$ws.Range("K55").Value = 481
$ws.Range("L55").Value = 2626.2
$ws.Range("M55").Value = -308
$ws.Range("N55").Value = -2972.2

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H132").Value = 3207.6316
$ws.Range("I132").Value = 2999.5334
$ws.Range("J132").Value = 3988
$ws.Range("K132").Value = 8998.600199999999
$ws.Range("L132").Value = 11964
$ws.Range("M132").Value = -6468.600199999999
$ws.Range("N132").Value = -17024

$ws.Range("H140").Value = 119495
$ws.Range("J140").Value = 119495
$ws.Range("L140").Value = 119495
$ws.Range("N140").Value = -129855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4457

$ws.Range("H65").Value = 4457

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 1837
$ws.Range("J107").Value = 1499
$ws.Range("L107").Value = 4497
$ws.Range("N107").Value = -8337

$ws.Range("H109").Value = 105000
$ws.Range("J109").Value = 105000
$ws.Range("L109").Value = 105000
$ws.Range("N109").Value = -107774

$ws.Range("H122").Value = 65379.695
$ws.Range("I122").Value = 78260.63
$ws.Range("K122").Value = 234781.89
$ws.Range("M122").Value = -232331.89

$ws.Range("H126").Value = 5033.9165
$ws.Range("I126").Value = 4271.143
$ws.Range("K126").Value = 12813.429
$ws.Range("M126").Value = -10343.429

$ws.Range("H136").Value = 17253.783
$ws.Range("I136").Value = 22177.74
$ws.Range("K136").Value = 66533.22
$ws.Range("M136").Value = -63983.22

Write-Host "Applied Sheets market data refresh"